# Update "Attributes" (column D) values to reflect per-action attribute
# similarity checks (set notation) instead of the whole-keyword-list
# comparison, and drop the extra (now-redundant) keyword from a couple of
# "Functionality" (column C) cells — per commit message:
# "checking similarity aganist each action(keyword) instead of comparing
# whole keywords"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Functionality column (C) simplifications
$ws.Range("C5").Value = "['Pay Bills ']"
$ws.Range("C7").Value = "['Cheque Services']"
$ws.Range("C9").Value = "['Limit Cash']"

# Attributes column (D) — now expressed as per-action sets
$ws.Range("D2").Value  = "{'Cust_Addr', 'Loan_Amt', 'Cus_Nme', 'Acc_num'}{'Debit_pin', 'Amt_avail', 'From_AcctNum', 'To_AcctNum', 'Acc_num'}"
$ws.Range("D3").Value  = "{'Debit_pin', 'Bill_type', 'Max_limit', 'Cus_Nme', 'Acc_num'}{'Cust_Addr', 'Loan_Amt', 'Cus_Nme', 'Acc_num'}"
$ws.Range("D4").Value  = "{'Amt_trnsfr', 'Amt_avail', 'Cus_Nme', 'To_AcctNum', 'From_AcctNum'}{'Cust_Addr', 'Loan_Amt', 'Cus_Nme', 'Acc_num'}"
$ws.Range("D5").Value  = "{'Debit_pin', 'From_AcctNum', 'Amt_avail', 'Bill_type', 'Acc_num'}"
$ws.Range("D6").Value  = "{'Debit_pin', 'From_AcctNum', 'Cus_Nme', 'To_AcctNum', 'Acc_num'}{'Loan_Amt', 'Loan_purp', 'Cred_Score'}"
$ws.Range("D7").Value  = "{'Cust_Addr', 'Loan_Amt', 'Cus_Nme', 'Acc_num'}"
$ws.Range("D8").Value  = "{'Max_limit', 'Acc_type', 'Cus_Nme', 'Acc_num'}{'Debit_pin', 'Acc_type', 'Cus_Nme', 'Amt_deposit', 'Acc_num'}"
$ws.Range("D9").Value  = "{'Debit_pin', 'Acc_type', 'Cus_Nme', 'Amt_deposit', 'Acc_num'}"
$ws.Range("D10").Value = "{'Acc_num'}{'Cust_Addr', 'Loan_Amt', 'Cus_Nme', 'Acc_num'}"
$ws.Range("D11").Value = "{'Loan_Amt', 'Loan_purp', 'Cred_Score'}{'Acc_num'},{'Debit_pin', 'From_AcctNum', 'Cus_Nme', 'To_AcctNum', 'Acc_num'}{'Loan_Amt', 'Loan_purp', 'Cred_Score'}"

# D12 was already empty (a stray typed-but-blank cell) — clear it outright.
$ws.Range("D12").ClearContents()

# Widen the Functionality column now that its values changed.
$ws.Columns.Item(3).ColumnWidth = 28
